$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 29
$ws.Range("H29").Value = 1000
$ws.Range("I29").Value = 1000
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 3000
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -2719
$ws.Range("N29").Value = ""
# Row 38
$ws.Range("H38").Value = 453.2
$ws.Range("I38").Value = 453.2
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 1359.6
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -987.5999999999999
$ws.Range("N38").Value = ""
# Row 40
$ws.Range("H40").Value = 5341.5835
$ws.Range("I40").Value = 3788.7778
$ws.Range("J40").Value = 10000
$ws.Range("K40").Value = 3788.7778
$ws.Range("L40").Value = 10000
$ws.Range("M40").Value = -3613.7778
# Row 87
$ws.Range("H87").Value = 69768
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 69768
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 69768
$ws.Range("N87").Value = -72264
# Row 90
$ws.Range("H90").Value = 69768
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 69768
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 209304
$ws.Range("N90").Value = -221784
# Row 98
$ws.Range("H98").Value = 2472.75
$ws.Range("I98").Value = 838.5
$ws.Range("J98").Value = 10644
$ws.Range("K98").Value = 838.5
$ws.Range("L98").Value = 10644
$ws.Range("M98").Value = 659.5
# Row 116
$ws.Range("H116").Value = 9900
$ws.Range("I116").Value = 9900
$ws.Range("J116").Value = 9900
$ws.Range("K116").Value = 9900
$ws.Range("L116").Value = 9900
$ws.Range("M116").Value = -6458
$ws.Range("N116").Value = -16784
# Row 122
$ws.Range("H122").Value = 2472.75
$ws.Range("I122").Value = 838.5
$ws.Range("J122").Value = 10644
$ws.Range("K122").Value = 2515.5
$ws.Range("L122").Value = 31932
$ws.Range("M122").Value = -65.5
# Row 132
$ws.Range("H132").Value = 2293.75
$ws.Range("I132").Value = 2293.75
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 6881.25
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -4351.25
# Row 135
$ws.Range("H135").Value = 958.4666999999999
$ws.Range("I135").Value = 678.6923
$ws.Range("J135").Value = 2777
$ws.Range("K135").Value = 6108.2307
$ws.Range("L135").Value = 24993
$ws.Range("M135").Value = -3573.2307
$ws.Range("N135").Value = -30063
# Row 138
$ws.Range("H138").Value = 6763.4194
$ws.Range("I138").Value = 4666.6665
$ws.Range("J138").Value = 6988.0713
$ws.Range("K138").Value = 13999.9995
$ws.Range("L138").Value = 20964.2139
$ws.Range("M138").Value = -8859.999500000002
$ws.Range("N138").Value = -31244.2139

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 9190.727999999999
$ws.Range("I32").Value = 9190.727999999999
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 9190.727999999999
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -8903.727999999999
# Row 74
$ws.Range("H74").Value = 24607
$ws.Range("I74").Value = 21004.666
$ws.Range("J74").Value = 30010.5
$ws.Range("K74").Value = 21004.666
$ws.Range("L74").Value = 30010.5
$ws.Range("M74").Value = -20130.666
# Row 77
$ws.Range("H77").Value = 24607
$ws.Range("I77").Value = 21004.666
$ws.Range("J77").Value = 30010.5
$ws.Range("K77").Value = 105023.33
$ws.Range("L77").Value = 150052.5
$ws.Range("M77").Value = -100655.33
# Row 102
$ws.Range("H102").Value = 3221.5
$ws.Range("I102").Value = 3221.5
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 3221.5
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -1599.5
$ws.Range("N102").Value = ""
# Row 132
$ws.Range("H132").Value = 3375
$ws.Range("I132").Value = 2747.5
$ws.Range("J132").Value = 4379
$ws.Range("K132").Value = 8242.5
$ws.Range("L132").Value = 13137
$ws.Range("M132").Value = -5712.5

$ws = $wb.Worksheets.Item("BSM")
# Row 7
$ws.Range("H7").Value = 16141.143
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 16141.143
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 16141.143
$ws.Range("M7").Value = ""
$ws.Range("N7").Value = -16367.143
# Row 20
$ws.Range("H20").Value = 8253
$ws.Range("I20").Value = 10966
$ws.Range("J20").Value = 7090.2856
$ws.Range("K20").Value = 10966
$ws.Range("L20").Value = 7090.2856
$ws.Range("M20").Value = -10719
# Row 134
$ws.Range("H134").Value = 1102.7142
$ws.Range("I134").Value = 1102.7142
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 3308.1426
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -773.1425999999997

$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 226.11111
$ws.Range("I7").Value = 160.5
$ws.Range("J7").Value = 278.6
$ws.Range("K7").Value = 160.5
$ws.Range("L7").Value = 278.6
$ws.Range("M7").Value = -47.5
$ws.Range("N7").Value = -504.6
# Row 31
$ws.Range("H31").Value = 2083.5
$ws.Range("I31").Value = 1681
$ws.Range("J31").Value = 3291
$ws.Range("K31").Value = 1681
$ws.Range("L31").Value = 3291
$ws.Range("M31").Value = -1386
$ws.Range("N31").Value = -3881
# Row 34
$ws.Range("H34").Value = 2083.5
$ws.Range("I34").Value = 1681
$ws.Range("J34").Value = 3291
$ws.Range("K34").Value = 1681
$ws.Range("L34").Value = 3291
$ws.Range("M34").Value = -1479
$ws.Range("N34").Value = -3695
# Row 59
$ws.Range("H59").Value = 23906.818
$ws.Range("I59").Value = 13000
$ws.Range("J59").Value = 24997.5
$ws.Range("K59").Value = 13000
$ws.Range("L59").Value = 24997.5
$ws.Range("M59").Value = -11855
$ws.Range("N59").Value = -27287.5
# Row 134
$ws.Range("H134").Value = 2495.6
$ws.Range("I134").Value = 2091
$ws.Range("J134").Value = 2958
$ws.Range("K134").Value = 6273
$ws.Range("L134").Value = 8874
$ws.Range("M134").Value = -3738

$ws = $wb.Worksheets.Item("CUL")
# Row 122
$ws.Range("H122").Value = 1794.4
$ws.Range("I122").Value = 994.25
$ws.Range("J122").Value = 2327.8333
$ws.Range("K122").Value = 8948.25
$ws.Range("L122").Value = 20950.4997
$ws.Range("M122").Value = -6498.25
$ws.Range("N122").Value = -25850.4997
# Row 129
$ws.Range("H129").Value = 1282.2222
$ws.Range("I129").Value = 974
$ws.Range("J129").Value = 1898.6666
$ws.Range("K129").Value = 2922
$ws.Range("L129").Value = 5695.9998
$ws.Range("M129").Value = 2078
$ws.Range("N129").Value = -15695.9998
# Row 131
$ws.Range("H131").Value = 4133
$ws.Range("I131").Value = 3995
$ws.Range("J131").Value = 4142.857
$ws.Range("K131").Value = 11985
$ws.Range("L131").Value = 12428.571
$ws.Range("M131").Value = -6945
$ws.Range("N131").Value = -22508.571

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 3563.8333
$ws.Range("I80").Value = 3577.4
$ws.Range("J80").Value = 3496
$ws.Range("K80").Value = 3577.4
$ws.Range("L80").Value = 3496
$ws.Range("M80").Value = -2579.4
# Row 83
$ws.Range("H83").Value = 3563.8333
$ws.Range("I83").Value = 3577.4
$ws.Range("J83").Value = 3496
$ws.Range("K83").Value = 17887
$ws.Range("L83").Value = 17480
$ws.Range("M83").Value = -12895
# Row 132
$ws.Range("H132").Value = 3828.8462
$ws.Range("I132").Value = 2398.2856
$ws.Range("J132").Value = 5497.8335
$ws.Range("K132").Value = 7194.8568
$ws.Range("L132").Value = 16493.5005
$ws.Range("M132").Value = -4664.8568
# Row 141
$ws.Range("H141").Value = 97999
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 97999
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 97999
$ws.Range("N141").Value = -108359

$ws = $wb.Worksheets.Item("LTW")
# Row 18
$ws.Range("H18").Value = 18000
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 18000
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 18000
$ws.Range("N18").Value = -18344
# Row 20
$ws.Range("H20").Value = 5000000
$ws.Range("I20").Value = 5000000
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 5000000
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -4999774
$ws.Range("N20").Value = ""
# Row 61
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").Value = ""
# Row 82
$ws.Range("H82").Value = 2275.7
$ws.Range("I82").Value = 2238.8
$ws.Range("J82").Value = 2312.6
$ws.Range("K82").Value = 2238.8
$ws.Range("L82").Value = 2312.6
$ws.Range("M82").Value = -1877.8
# Row 85
$ws.Range("H85").Value = 2275.7
$ws.Range("I85").Value = 2238.8
$ws.Range("J85").Value = 2312.6
$ws.Range("K85").Value = 2238.8
$ws.Range("L85").Value = 2312.6
$ws.Range("M85").Value = -990.8000000000002
# Row 113
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").Value = ""
# Row 122
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").Value = ""
# Row 136
$ws.Range("H136").Value = 5199.5
$ws.Range("I136").Value = 1918.4
$ws.Range("J136").Value = 10668
$ws.Range("K136").Value = 5755.200000000001
$ws.Range("L136").Value = 32004
$ws.Range("M136").Value = -3205.200000000001
$ws.Range("N136").Value = -37104

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 7442.875
$ws.Range("I62").Value = 2950
$ws.Range("J62").Value = 8940.5
$ws.Range("K62").Value = 2950
$ws.Range("L62").Value = 8940.5
$ws.Range("M62").Value = -2326
$ws.Range("N62").Value = -10188.5
# Row 65
$ws.Range("H65").Value = 7442.875
$ws.Range("I65").Value = 2950
$ws.Range("J65").Value = 8940.5
$ws.Range("K65").Value = 14750
$ws.Range("L65").Value = 44702.5
$ws.Range("M65").Value = -11630
$ws.Range("N65").Value = -50942.5
# Row 81
$ws.Range("H81").Value = 5071.273
$ws.Range("I81").Value = 2826.2856
$ws.Range("J81").Value = 9000
$ws.Range("K81").Value = 5652.5712
$ws.Range("L81").Value = 18000
$ws.Range("M81").Value = -4591.5712
$ws.Range("N81").Value = -20122
# Row 84
$ws.Range("H84").Value = 5071.273
$ws.Range("I84").Value = 2826.2856
$ws.Range("J84").Value = 9000
$ws.Range("K84").Value = 28262.856
$ws.Range("L84").Value = 90000
$ws.Range("M84").Value = -22958.856
$ws.Range("N84").Value = -100608
# Row 140
$ws.Range("H140").Value = 100000
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 100000
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 100000
$ws.Range("N140").Value = -110360
